$wb = $excel.ActiveWorkbook

# --- "Chris K." sheet: remove the old "Brush" row (row 4), shifting the
#     "Oatmeal Soap" row up, then update the remaining rows/vendor cells ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(4).Delete()
$ws2.Range("A3").Value = "Brush"
$ws2.Range("C3").Value = 4.99
$ws2.Range("E2").Value = "Premium Cat Food"
$ws2.Range("I2").Value = 21.98
$ws2.Range("K2").Value = 42.18

# --- "Mary M." sheet: remove the "Fruits" and "Vegetables" rows (old rows
#     3 and 4), shifting "Basic Dog Food" up to row 3, then update vendor
#     cells ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(3).Delete()
$ws3.Rows.Item(3).Delete()
$ws3.Range("A3").Value = "Basic Dog Food"
$ws3.Range("C3").Value = 8.99
$ws3.Range("E2").Value = "Fruits"
$ws3.Range("E3").Value = "Vegetables"
$ws3.Range("I2").Value = 20.98
$ws3.Range("K2").Value = 40.98
